# Fill in the single data row of the "transaccion" dataset spec sheet.
# Row 1 (headers) is already populated; row 2 currently holds empty
# placeholder cells that need real sample values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:I2")

# Several of the values look numeric/date-like ("896655273", "3568999",
# "150000", "2019-05-31"). Force the range to Text format first so Excel
# stores them as literal strings instead of silently coercing them to
# numbers/dates.
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "896655273"
$ws.Range("B2").Value = "3568999-1"
$ws.Range("C2").Value = "3568999"
$ws.Range("D2").Value = "Prestamo BID AR1505"
$ws.Range("E2").Value = "2019-05-31"
$ws.Range("F2").Value = "ARS"
$ws.Range("G2").Value = "150000"
$ws.Range("H2").Value = "Ministerio de modernización"
$ws.Range("I2").Value = "HAL2000"

# Restore the default "Normal" cell style so the row doesn't pick up a
# lingering Text number-format style (matches the source file, where the
# row has no explicit style applied).
$dataRange.Style = "Normal"
